$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert two new columns (K:L) for the new "quantityValueDropDown" and
#    "frequentlyOrderedItem" fields. This shifts the existing rA*/product*/
#    orderTime columns (old K..Q) two slots to the right (new M..S), and
#    correctly carries along the existing K2/L2/K3/L3 placeholder cells,
#    hyperlinks and styling.
# ---------------------------------------------------------------------------
$ws.Columns("K:L").Insert()

# New header cells for the inserted columns.
$ws.Range("K1").Value = "quantityValueDropDown"
$ws.Range("L1").Value = "frequentlyOrderedItem"

# Give the two new columns a sensible best-fit-like width (best effort).
$ws.Columns("K").ColumnWidth = 22.67
$ws.Columns("L").ColumnWidth = 26.67

# ---------------------------------------------------------------------------
# 2. Existing row 4 ("TC 3") gets a corrected productPrize value (now in the
#    shifted column R) -- was 5.40, becomes 15.60. The leading apostrophe
#    keeps it stored as text (matching the other price cells) instead of a
#    number, while preserving the existing quote-prefix cell style.
# ---------------------------------------------------------------------------
$ws.Range("R4").Value = "'15.60"

# ---------------------------------------------------------------------------
# 3. Brand-new test case rows 5 and 6.
#    Style "1" (quote-prefixed General) is applied simply by prefixing the
#    literal text with an apostrophe - the engine reuses the existing
#    quote-prefix style automatically. Style "2" is the built-in "Hyperlink"
#    named style, applied explicitly before setting the value.
# ---------------------------------------------------------------------------

# --- Row 5 ("TC 4") ---
$ws.Range("A5").Value = "'4"
$ws.Range("B5").Style = "Hyperlink"
$ws.Range("B5").Value = "uber+test+developer@qikserve.com"
$ws.Range("C5").Value = "'cT;p3kx_"
$ws.Range("D5").Value = "'Matt"
$ws.Range("E5").Value = "'2 Leman Street"
$ws.Range("F5").Value = "'Coffee Frappuccino®"
$ws.Range("G5").Value = "Grande"
$ws.Range("H5").Value = "Semi-Skimmed Milk"
$ws.Range("K5").Value = "'3"
$ws.Range("L5").Value = "'"
$ws.Range("L5").ClearContents()
$ws.Range("M5").Value = "SSELVARA"
$ws.Range("N5").Value = "SBU"
$ws.Range("O5").Style = "Hyperlink"
$ws.Range("O5").Value = "Qwerty@6098"
$ws.Range("P5").Value = "West Bromwich Drive"
$ws.Range("Q5").Value = "Caramel Frappucino® Blended Beverage"
$ws.Range("R5").Value = "'15.60"
$ws.Range("S5").Value = "'14:16"

# --- Row 6 ("TC 5") ---
$ws.Range("A6").Value = "'5"
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("B6").Value = "uber+test+developer@qikserve.com"
$ws.Range("C6").Value = "'cT;p3kx_"
$ws.Range("D6").Value = "'Matt"
$ws.Range("E6").Value = "'2 Leman Street"
$ws.Range("F6").Value = "Filter Coffee"
$ws.Range("G6").Value = "Grande"
$ws.Range("L6").Value = "Smoked Bacon Roll"
$ws.Range("M6").Value = "SSELVARA"
$ws.Range("N6").Value = "SBU"
$ws.Range("O6").Style = "Hyperlink"
$ws.Range("O6").Value = "Qwerty@6098"
$ws.Range("P6").Value = "West Bromwich Drive"
$ws.Range("Q6").Value = "Filter Coffee"
$ws.Range("R6").Value = "'7.50"
$ws.Range("S6").Value = "'11:55"

# ---------------------------------------------------------------------------
# 4. Hyperlinks for the new rows (mirroring the B2/O2, B3/O3, B4/O4 pattern).
#    Hyperlinks.Add always (re)applies its own font variant of the
#    "Hyperlink" style, so re-assert the canonical named style afterwards to
#    keep every hyperlinked cell sharing the same style index as B2/O2/etc.
# ---------------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("B5"), "mailto:uber+test+developer@qikserve.com")
$ws.Hyperlinks.Add($ws.Range("O5"), "mailto:Qwerty@6098")
$ws.Hyperlinks.Add($ws.Range("B6"), "mailto:uber+test+developer@qikserve.com")
$ws.Hyperlinks.Add($ws.Range("O6"), "mailto:Qwerty@6098")

$ws.Range("B5").Style = "Hyperlink"
$ws.Range("O5").Style = "Hyperlink"
$ws.Range("B6").Style = "Hyperlink"
$ws.Range("O6").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# 5. Selection / view bookkeeping to mirror the saved workbook state.
# ---------------------------------------------------------------------------
$ws.Range("U22").Select()
